{"js": "const body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraphs we need to touch by their current text.\nlet nameIdx = -1, sahilIdx = -1, descIdx = -1, bibIdx = -1, appIdx = -1;\nfor (let i = 0; i < paras.items.length; i++) {\n  const t = paras.items[i].text;\n  if (t === \"ARYYAMA KUMAR JANA\") nameIdx = i;\n  else if (t === \"Sahil Santosh Patil (He/Him)\") sahilIdx = i;\n  else if (t === \"Description of Tests: describe the tests that were performed\") descIdx = i;\n  else if (t === \"Bibliography: any resources(a lot of them) that were used\") bibIdx = i;\n  else if (t === \"Appendix(typescript): the garbage generated by the test\") appIdx = i;\n}\n\n// 1) \"ARYYAMA KUMAR JANA\" -> \"Aryyama Kumar Jana\"\nif (nameIdx !== -1) {\n  paras.items[nameIdx].insertText(\"Aryyama Kumar Jana\", \"Replace\");\n}\n\n// 2) \"Sahil Santosh Patil (He/Him)\" -> \"Sahil Santosh Patil\" (drop the \"(He/Him)\")\nif (sahilIdx !== -1) {\n  paras.items[sahilIdx].insertText(\"Sahil Santosh Patil\", \"Replace\");\n}\n\n// 3) Split \"Description of Tests: describe the tests that were performed\"\n//    into \"Description of Tests: \" followed by a new paragraph describing\n//    Test 1 and Test 2.\nif (descIdx !== -1) {\n  const descPara = paras.items[descIdx];\n  descPara.clear();\n  descPara.insertText(\"Description of Tests: \", \"Replace\");\n  const newPara = descPara.insertParagraph(\"\", \"After\");\n  newPara.insertText(\n    \"\\tTest 1 tests basic buffer management operations, specifically valid and legal operations, by allocating new pages and writing to those pages. The pages are then read back before being freed. This tests functionality of the database buffer manager such as its replacement algorithm and page management. Test 2 also tests the buffer manager by testing its expected failures. These include things such as pinning a greater number of pages than frames, freeing doubly pinned pages, and unpinning pages not currently in the buffer pool.\",\n    \"Replace\"\n  );\n}\n\n// 4) Bibliography / Appendix keep the same wording, but the stray\n//    mid-sentence proofing marks (gramStart/gramEnd) that split the text\n//    across extra runs are gone in the final version, so rewrite each\n//    paragraph as a single clean run with identical text.\nif (bibIdx !== -1) {\n  const bibPara = paras.items[bibIdx];\n  bibPara.clear();\n  bibPara.insertText(\"Bibliography: any resources(a lot of them) that were used\", \"Replace\");\n}\nif (appIdx !== -1) {\n  const appPara = paras.items[appIdx];\n  appPara.clear();\n  appPara.insertText(\"Appendix(typescript): the garbage generated by the test\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\nfunction FindParaIndexByText($doc, $text) {\n    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {\n        $t = $doc.Paragraphs($i).Range.Text.TrimEnd(\"`r\")\n        if ($t -eq $text) { return $i }\n    }\n    return -1\n}\n\n# Replace a paragraph's text in place (keeps the paragraph mark / pPr),\n# without leaving stray proofErr siblings behind.\nfunction ReplaceParaTextClean($doc, $paraIndex, $newText) {\n    $oldPara = $doc.Paragraphs($paraIndex)\n    $oldPara.Range.InsertParagraphAfter()\n    $newPara = $doc.Paragraphs($paraIndex + 1)\n    $newPara.Range.Text = $newText\n    $oldPara.Range.Delete()\n}\n\n# 1) \"ARYYAMA KUMAR JANA\" -> \"Aryyama Kumar Jana\"\n$idx = FindParaIndexByText $d \"ARYYAMA KUMAR JANA\"\nif ($idx -gt 0) {\n    $r = $d.Paragraphs($idx).Range\n    $r.MoveEnd(1, -1) | Out-Null\n    $r.Text = \"Aryyama Kumar Jana\"\n}\n\n# 2) \"Sahil Santosh Patil (He/Him)\" -> \"Sahil Santosh Patil\"\n$idx = FindParaIndexByText $d \"Sahil Santosh Patil (He/Him)\"\nif ($idx -gt 0) {\n    $r = $d.Paragraphs($idx).Range\n    $r.MoveEnd(1, -1) | Out-Null\n    $r.Text = \"Sahil Santosh Patil\"\n}\n\n# 3) Split \"Description of Tests: describe the tests that were performed\"\n#    into \"Description of Tests: \" followed by a new paragraph describing\n#    Test 1 and Test 2.\n$idx = FindParaIndexByText $d \"Description of Tests: describe the tests that were performed\"\nif ($idx -gt 0) {\n    $testsText = \"Test 1 tests basic buffer management operations, specifically valid and legal operations, by allocating new pages and writing to those pages. The pages are then read back before being freed. This tests functionality of the database buffer manager such as its replacement algorithm and page management. Test 2 also tests the buffer manager by testing its expected failures. These include things such as pinning a greater number of pages than frames, freeing doubly pinned pages, and unpinning pages not currently in the buffer pool.\"\n\n    $oldPara = $d.Paragraphs($idx)\n    $oldPara.Range.InsertParagraphAfter()\n    $firstNew = $d.Paragraphs($idx + 1)\n    $firstNew.Range.Text = \"Description of Tests: \"\n    $firstNew.Range.InsertParagraphAfter()\n    $secondNew = $d.Paragraphs($idx + 2)\n    $secondNew.Range.Text = \"`t\" + $testsText\n    $oldPara.Range.Delete()\n}\n\n# 4) Bibliography / Appendix keep the same wording, but the stray\n#    mid-sentence proofing marks (gramStart/gramEnd) that split the text\n#    across extra runs are gone in the final version, so rewrite each\n#    paragraph as a single clean run with identical text.\n$idx = FindParaIndexByText $d \"Bibliography: any resources(a lot of them) that were used\"\nif ($idx -gt 0) {\n    ReplaceParaTextClean $d $idx \"Bibliography: any resources(a lot of them) that were used\"\n}\n\n$idx = FindParaIndexByText $d \"Appendix(typescript): the garbage generated by the test\"\nif ($idx -gt 0) {\n    ReplaceParaTextClean $d $idx \"Appendix(typescript): the garbage generated by the test\"\n}\n"}
